# The workbook originally has a generically-named "Sheet4" tab holding the
# maintenance-category lookup data (id / Description rows for YA02, YBA1,
# YBA2 ...). Give it a descriptive name, matching the rest of the master
# data tabs (Equipment, Functional Location, Company Code, Cost Center,
# Plant, ...).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")
$ws.Name = "Maintenance Category"

# Bring the freshly-renamed sheet into focus, which is what actually moves
# the workbook's "active tab" / selected-tab state onto it.
$ws.Activate()
